# The sheet used to store a hand-typed "id" column (Sean001, Sean002, ...)
# in column A and the "author" column (always "Sean cheng") last, in
# column H. Going forward ids are generated programmatically
# (yyyymmddhhmmss + index) instead of being stored in the sheet, so the
# old "id" column goes away and "author" becomes the new first column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old hand-typed "id" column (A) ...
$ws.Columns("A").Delete()

# ... then move the "author" column (now the last column, G, after the
# shift above) to the front of the sheet, so it becomes the new column A.
$ws.Columns("G").Cut()
$ws.Columns("A").Insert()
